$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 239.44444
$ws.Range("I2").Value = 432
$ws.Range("J2").Value = 85.40000000000001
$ws.Range("K2").Value = 432
$ws.Range("L2").Value = 85.40000000000001
$ws.Range("M2").Value = -319
$ws.Range("N2").Value = -311.4
$ws.Range("H16").Value = 9066.666999999999
$ws.Range("J16").Value = 25000
$ws.Range("L16").Value = 25000
$ws.Range("N16").Value = -25460
$ws.Range("H17").Value = 1843.4
$ws.Range("J17").Value = 1919.5714
$ws.Range("L17").Value = 5758.7142
$ws.Range("N17").Value = -6094.7142
$ws.Range("H18").Value = 724.63635
$ws.Range("I18").Value = 351
$ws.Range("J18").Value = 2406
$ws.Range("K18").Value = 351
$ws.Range("L18").Value = 2406
$ws.Range("M18").Value = -67
$ws.Range("N18").Value = -2974
$ws.Range("H34").Value = 16232.75
$ws.Range("I34").Value = 16232.75
$ws.Range("K34").Value = 16232.75
$ws.Range("M34").Value = -16029.75
$ws.Range("H36").Value = 16232.75
$ws.Range("I36").Value = 16232.75
$ws.Range("K36").Value = 16232.75
$ws.Range("M36").Value = -15517.75
$ws.Range("H38").Value = 5963.483
$ws.Range("J38").Value = 10218.875
$ws.Range("L38").Value = 30656.625
$ws.Range("N38").Value = -31400.625
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 7991.885
$ws.Range("I64").Value = 7090.5713
$ws.Range("J64").Value = 8323.947
$ws.Range("K64").Value = 7090.5713
$ws.Range("L64").Value = 8323.947
$ws.Range("M64").Value = -6842.5713
$ws.Range("N64").Value = -8819.947
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 7991.885
$ws.Range("I67").Value = 7090.5713
$ws.Range("J67").Value = 8323.947
$ws.Range("K67").Value = 7090.5713
$ws.Range("L67").Value = 8323.947
$ws.Range("M67").Value = -6232.5713
$ws.Range("N67").Value = -10039.947
$ws.Range("H92").Value = 2976813
$ws.Range("I92").Value = 391.6842
$ws.Range("J92").Value = 14287213
$ws.Range("K92").Value = 391.6842
$ws.Range("L92").Value = 14287213
$ws.Range("M92").Value = 856.3158000000001
$ws.Range("N92").Value = -14289709
$ws.Range("H106").Value = 62506.57
$ws.Range("J106").Value = 26497.5
$ws.Range("L106").Value = 26497.5
$ws.Range("N106").Value = -27759.5
$ws.Range("H107").Value = 2338.8333
$ws.Range("I107").Value = 2298
$ws.Range("K107").Value = 2298
$ws.Range("M107").Value = -378
$ws.Range("H111").Value = 1924.3077
$ws.Range("I111").Value = 1375.7778
$ws.Range("K111").Value = 4127.3334
$ws.Range("M111").Value = -1060.3334
$ws.Range("H112").Value = 1258.1072
$ws.Range("J112").Value = 1258.1072
$ws.Range("L112").Value = 3774.3216
$ws.Range("N112").Value = -5990.321599999999
$ws.Range("H115").Value = 3250.5715
$ws.Range("I115").Value = 3250.5715
$ws.Range("K115").Value = 9751.7145
$ws.Range("M115").Value = -8184.7145
$ws.Range("H125").Value = 1494.3334
$ws.Range("I125").Value = 1492.5
$ws.Range("K125").Value = 13432.5
$ws.Range("M125").Value = -10972.5
$ws.Range("H132").Value = 3087.3572
$ws.Range("I132").Value = 3059.4614
$ws.Range("K132").Value = 9178.3842
$ws.Range("M132").Value = -6648.3842
$ws.Range("H138").Value = 90912680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2315.5942
$ws.Range("I32").Value = 1463.8871
$ws.Range("K32").Value = 1463.8871
$ws.Range("M32").Value = -1176.8871
$ws.Range("H45").Value = 4415.385
$ws.Range("I45").Value = 3903.1667
$ws.Range("J45").Value = 4854.4287
$ws.Range("K45").Value = 3903.1667
$ws.Range("L45").Value = 4854.4287
$ws.Range("M45").Value = -3526.1667
$ws.Range("N45").Value = -5608.4287
$ws.Range("H63").Value = 1883
$ws.Range("I63").Value = 1918.9231
$ws.Range("K63").Value = 1918.9231
$ws.Range("M63").Value = -1232.9231
$ws.Range("H66").Value = 1883
$ws.Range("I66").Value = 1918.9231
$ws.Range("K66").Value = 9594.6155
$ws.Range("M66").Value = -6162.6155
$ws.Range("H74").Value = 2912.9443
$ws.Range("I74").Value = 2347.6296
$ws.Range("J74").Value = 4608.8887
$ws.Range("K74").Value = 2347.6296
$ws.Range("L74").Value = 4608.8887
$ws.Range("M74").Value = -1473.6296
$ws.Range("N74").Value = -6356.8887
$ws.Range("H77").Value = 2912.9443
$ws.Range("I77").Value = 2347.6296
$ws.Range("J77").Value = 4608.8887
$ws.Range("K77").Value = 11738.148
$ws.Range("L77").Value = 23044.4435
$ws.Range("M77").Value = -7370.148000000001
$ws.Range("N77").Value = -31780.4435
$ws.Range("H97").Value = 402.08334
$ws.Range("I97").Value = 331.375
$ws.Range("J97").Value = 543.5
$ws.Range("K97").Value = 331.375
$ws.Range("L97").Value = 543.5
$ws.Range("M97").Value = 164.625
$ws.Range("N97").Value = -1535.5
$ws.Range("H102").Value = 21319.76
$ws.Range("I102").Value = 1408.909
$ws.Range("K102").Value = 1408.909
$ws.Range("M102").Value = 213.0909999999999
$ws.Range("H110").Value = 1085.9642
$ws.Range("I110").Value = 974.04346
$ws.Range("J110").Value = 1600.8
$ws.Range("K110").Value = 974.04346
$ws.Range("L110").Value = 1600.8
$ws.Range("M110").Value = 1070.95654
$ws.Range("N110").Value = -5690.8
$ws.Range("H122").Value = 4498.731
$ws.Range("I122").Value = 4077.1904
$ws.Range("K122").Value = 12231.5712
$ws.Range("M122").Value = -9781.5712
$ws.Range("H132").Value = 10219.588
$ws.Range("I132").Value = 6305.2896
$ws.Range("J132").Value = 21661.385
$ws.Range("K132").Value = 18915.8688
$ws.Range("L132").Value = 64984.155
$ws.Range("M132").Value = -16385.8688
$ws.Range("N132").Value = -70044.155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 683.05884
$ws.Range("J80").Value = 607.1818
$ws.Range("L80").Value = 607.1818
$ws.Range("N80").Value = -2603.1818
$ws.Range("H83").Value = 683.05884
$ws.Range("J83").Value = 607.1818
$ws.Range("L83").Value = 3035.909
$ws.Range("N83").Value = -13019.909
$ws.Range("H86").Value = 3120.6316
$ws.Range("I86").Value = 2784.077
$ws.Range("K86").Value = 2784.077
$ws.Range("M86").Value = -1661.077
$ws.Range("H89").Value = 3120.6316
$ws.Range("I89").Value = 2784.077
$ws.Range("K89").Value = 13920.385
$ws.Range("M89").Value = -8304.385000000002
$ws.Range("H99").Value = 100001830
$ws.Range("I99").Value = 200001000
$ws.Range("K99").Value = 200001000
$ws.Range("M99").Value = -199999502
$ws.Range("H105").Value = 3141.258
$ws.Range("I105").Value = 3115
$ws.Range("K105").Value = 3115
$ws.Range("M105").Value = -1368
$ws.Range("H134").Value = 5161.4546
$ws.Range("I134").Value = 8999
$ws.Range("J134").Value = 4777.7
$ws.Range("K134").Value = 26997
$ws.Range("L134").Value = 14333.1
$ws.Range("M134").Value = -24462
$ws.Range("N134").Value = -19403.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4768.8184
$ws.Range("I99").Value = 1691
$ws.Range("J99").Value = 6527.5713
$ws.Range("K99").Value = 1691
$ws.Range("L99").Value = 6527.5713
$ws.Range("M99").Value = -193
$ws.Range("N99").Value = -9523.5713
$ws.Range("H122").Value = 2722.762
$ws.Range("J122").Value = 3042.3333
$ws.Range("L122").Value = 9126.999899999999
$ws.Range("N122").Value = -14026.9999
$ws.Range("H126").Value = 4768.8184
$ws.Range("I126").Value = 1691
$ws.Range("J126").Value = 6527.5713
$ws.Range("K126").Value = 5073
$ws.Range("L126").Value = 19582.7139
$ws.Range("M126").Value = -2603
$ws.Range("N126").Value = -24522.7139
$ws.Range("H132").Value = 3097.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3097.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9292.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -14352.5
$ws.Range("H134").Value = 1919.7858
$ws.Range("I134").Value = 1936.0834
$ws.Range("K134").Value = 5808.2502
$ws.Range("M134").Value = -3273.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 62500172
$ws.Range("I4").Value = 66666850
$ws.Range("K4").Value = 200000550
$ws.Range("M4").Value = -200000438
$ws.Range("H6").Value = 1617.6666
$ws.Range("I6").Value = 156.88889
$ws.Range("J6").Value = 6000
$ws.Range("K6").Value = 470.66667
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = -357.66667
$ws.Range("N6").Value = -18226
$ws.Range("H39").Value = 9947.200000000001
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H110").Value = 724.5
$ws.Range("I110").Value = 724.5
$ws.Range("K110").Value = 2173.5
$ws.Range("M110").Value = 1916.5
$ws.Range("H125").Value = 3399.25
$ws.Range("I125").Value = 3399.25
$ws.Range("K125").Value = 10197.75
$ws.Range("M125").Value = -5277.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 48722.715
$ws.Range("I70").Value = 56557.184
$ws.Range("J70").Value = 19996.334
$ws.Range("K70").Value = 56557.184
$ws.Range("L70").Value = 19996.334
$ws.Range("M70").Value = -56287.184
$ws.Range("N70").Value = -20536.334
$ws.Range("H73").Value = 48722.715
$ws.Range("I73").Value = 56557.184
$ws.Range("J73").Value = 19996.334
$ws.Range("K73").Value = 56557.184
$ws.Range("L73").Value = 19996.334
$ws.Range("M73").Value = -55621.184
$ws.Range("N73").Value = -21868.334
$ws.Range("H97").Value = 104
$ws.Range("I97").Value = 104
$ws.Range("K97").Value = 104
$ws.Range("M97").Value = 392
$ws.Range("H102").Value = 2376.1667
$ws.Range("I102").Value = 1717.0555
$ws.Range("K102").Value = 1717.0555
$ws.Range("M102").Value = -95.05549999999994
$ws.Range("H122").Value = 3098.0715
$ws.Range("I122").Value = 1669.6666
$ws.Range("J122").Value = 3487.6365
$ws.Range("K122").Value = 5008.9998
$ws.Range("L122").Value = 10462.9095
$ws.Range("M122").Value = -2558.9998
$ws.Range("N122").Value = -15362.9095
$ws.Range("H132").Value = 4398
$ws.Range("I132").Value = 3656
$ws.Range("K132").Value = 10968
$ws.Range("M132").Value = -8438

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3012.4285
$ws.Range("I22").Value = 1009.5
$ws.Range("J22").Value = 5683
$ws.Range("K22").Value = 1009.5
$ws.Range("L22").Value = 5683
$ws.Range("M22").Value = -714.5
$ws.Range("N22").Value = -6273
$ws.Range("H27").Value = 3012.4285
$ws.Range("I27").Value = 1009.5
$ws.Range("J27").Value = 5683
$ws.Range("K27").Value = 1009.5
$ws.Range("L27").Value = 5683
$ws.Range("M27").Value = -902.5
$ws.Range("N27").Value = -5897
$ws.Range("H31").Value = 2912.2666
$ws.Range("I31").Value = 2390.3333
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 2390.3333
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -2142.3333
$ws.Range("N31").Value = -5496
$ws.Range("H43").Value = 21500
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 21500
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 21500
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -21886
$ws.Range("H46").Value = 1834.5
$ws.Range("I46").Value = 1362.6666
$ws.Range("K46").Value = 1362.6666
$ws.Range("M46").Value = -1174.6666
$ws.Range("H132").Value = 3742
$ws.Range("I132").Value = 3336.0852
$ws.Range("K132").Value = 10008.2556
$ws.Range("M132").Value = -7478.2556
$ws.Range("H136").Value = 2631.8572
$ws.Range("I136").Value = 2225.375
$ws.Range("J136").Value = 6967.6665
$ws.Range("K136").Value = 6676.125
$ws.Range("L136").Value = 20902.9995
$ws.Range("M136").Value = -4126.125
$ws.Range("N136").Value = -26002.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33752.223
$ws.Range("I54").Value = 11885
$ws.Range("K54").Value = 11885
$ws.Range("M54").Value = -11365
$ws.Range("H62").Value = 19669.666
$ws.Range("J62").Value = 20799.5
$ws.Range("L62").Value = 20799.5
$ws.Range("N62").Value = -22047.5
$ws.Range("H65").Value = 19669.666
$ws.Range("J65").Value = 20799.5
$ws.Range("L65").Value = 103997.5
$ws.Range("N65").Value = -110237.5
$ws.Range("H81").Value = 5440.625
$ws.Range("I81").Value = 3272.5833
$ws.Range("J81").Value = 11944.75
$ws.Range("K81").Value = 6545.1666
$ws.Range("L81").Value = 23889.5
$ws.Range("M81").Value = -5484.1666
$ws.Range("N81").Value = -26011.5
$ws.Range("H84").Value = 5440.625
$ws.Range("I84").Value = 3272.5833
$ws.Range("J84").Value = 11944.75
$ws.Range("K84").Value = 32725.833
$ws.Range("L84").Value = 119447.5
$ws.Range("M84").Value = -27421.833
$ws.Range("N84").Value = -130055.5
$ws.Range("H107").Value = 790.73334
$ws.Range("I107").Value = 448.375
$ws.Range("K107").Value = 1345.125
$ws.Range("M107").Value = 574.875
$ws.Range("H113").Value = 1987.1111
$ws.Range("I113").Value = 1059.8
$ws.Range("J113").Value = 3146.25
$ws.Range("K113").Value = 3179.4
$ws.Range("L113").Value = 9438.75
$ws.Range("M113").Value = -1009.4
$ws.Range("N113").Value = -13778.75
$ws.Range("H126").Value = 15153062
$ws.Range("I126").Value = 19609502
$ws.Range("J126").Value = 1169.8
$ws.Range("K126").Value = 58828506
$ws.Range("L126").Value = 3509.4
$ws.Range("M126").Value = -58826036
$ws.Range("N126").Value = -8449.4
$ws.Range("H135").Value = 9999
$ws.Range("J135").Value = 9999
$ws.Range("L135").Value = 9999
$ws.Range("N135").Value = -20139
$ws.Range("H136").Value = 3918.8572
$ws.Range("I136").Value = 3551.6667
$ws.Range("J136").Value = 5265.222
$ws.Range("K136").Value = 10655.0001
$ws.Range("L136").Value = 15795.666
$ws.Range("M136").Value = -8105.000100000001
$ws.Range("N136").Value = -20895.666
